$d = $word.ActiveDocument

# --- 1. Add bookmark "_GoBack" (id 0) at the start of the second paragraph,
#        right before "You will be experimenting..." ---
$r1 = $d.Content
$find1 = $r1.Find
$find1.Execute("You will be experimenting", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Collapse(1)   # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $r1) | Out-Null

# --- 2. Split the run "... Labs are marked as ..." so that bookmark end
#        for _GoBack sits right after "Labs " (before "are marked") ---
$r2 = $d.Content
$find2 = $r2.Find
$find2.Execute("by completing the exercises below. Labs ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)   # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $r2) | Out-Null

# --- 3. Remove the old bookmarkStart/bookmarkEnd pair for "_GoBack"
#        that used to sit after "APIARY and " ---
$bm = $d.Bookmarks("_GoBack")
Write-Output "bookmark exists before cleanup: $($d.Bookmarks.Exists('_GoBack'))"
